# Adding CRM accuracy data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix casing of the existing "N/A" note to "n/a"
$ws.Range("B148").Value = "n/a"

# Correct the mislabeled CRM-opened batch note on row 149
$ws.Range("F149").Value = "CRM196_opened20240828"

# Widen column B to fit the data
$ws.Columns.Item(2).ColumnWidth = 16

# Append the new CRM accuracy data row
$ws.Range("A150").Value = 20241001
$ws.Range("B150").Value = 2207.22277248604
$ws.Range("C150").Value = 2215.3200000000002
$ws.Range("D150").Formula = "=100*(B150-C150)/C150"
$ws.Range("E150").Value = 196
$ws.Range("F150").Value = "CRM196_opened20240828"
